$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.979.70'
$ws.Range("E2").Value = '  -3.59%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.792.92'
$ws.Range("E3").Value = '  +0.72%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.72'
$ws.Range("E5").Value = '  -4.19%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.35'
$ws.Range("E6").Value = '  -5.42%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.793.16'
$ws.Range("E7").Value = '  +0.84%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  -1.41%  '

$ws.Range("E10").Value = '  -4.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.27'
$ws.Range("E11").Value = '  -1.74%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.464'
$ws.Range("E12").Value = '  -3.73%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.14'
$ws.Range("E13").Value = '  -5.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000244'
$ws.Range("E14").Value = '  -4.69%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.427.75'
$ws.Range("E15").Value = '  +0.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.790.97'
$ws.Range("E16").Value = '  +0.71%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.058.96'
$ws.Range("E17").Value = '  -3.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.15'
$ws.Range("E19").Value = '  -5.49%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.04'
$ws.Range("E20").Value = '  -3.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '488.56'
$ws.Range("E21").Value = '  -3.69%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.25'
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.731'
$ws.Range("E23").Value = '  +1.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.70'
$ws.Range("E24").Value = '  -2.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  -9.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000137'
$ws.Range("E26").Value = '  +1.68%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.20'
$ws.Range("E27").Value = '  -7.12%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.19'
$ws.Range("E28").Value = '  -10.56%  '

$ws.Range("E29").Value = '  -0.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.92'
$ws.Range("E30").Value = '  -0.57%  '

$ws.Range("E31").Value = '  -2.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.72'
$ws.Range("E32").Value = '  +6.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.70'
$ws.Range("E33").Value = '  -3.06%  '

$ws.Range("E34").Value = '  -4.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.08%  '

$ws.Range("E36").Value = '  -5.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.136'
$ws.Range("E37").Value = '  -3.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.76'
$ws.Range("E38").Value = '  -6.29%  '

$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.325'
$ws.Range("E39").Value = '  -8.49%  '

$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '453.68'
$ws.Range("E40").Value = '  +3.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '48.91'
$ws.Range("E41").Value = '  -2.48%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.99'
$ws.Range("E42").Value = '  -4.77%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.88'
$ws.Range("E43").Value = '  -7.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.25'
$ws.Range("E44").Value = '  -4.60%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.45'
$ws.Range("E45").Value = '  -9.24%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.832.23'
$ws.Range("E46").Value = '  -4.30%  '

$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '139.24'
$ws.Range("E47").Value = '  +1.02%  '

$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.02%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0351'
$ws.Range("E49").Value = '  -3.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.25'
$ws.Range("E50").Value = '  -4.55%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.22'
$ws.Range("E51").Value = '  +7.83%  '
